$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4: new data row
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 40804

# Row 11-12: new block (write in this order so shared-string indices line up)
$ws.Range("B11").Value = "CDC Patch to be released:"
$ws.Range("D11").Value = "EKA_CDC_1.0_Objects _0008.sql"
$ws.Range("D12").Value = "EKA_CDC_1.0_Objects _0007.sql"

$ws.Range("I4").Value = "EKA_METALS_PATCH_0145.txt"

# Selection update
$ws.Range("I4").Select()
